$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-7 from 45224 to 45233
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 3).Value = 45233
}
